# Update gh-pages to output generated at 456a3b4
# Updates the "想去人数" (number of people interested) column (F) values
# on the 展览, 演出, and 全部类型 sheets.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 156
$ws1.Range("F3").Value = 469
$ws1.Range("F4").Value = 11
$ws1.Range("F6").Value = 13
$ws1.Range("F8").Value = 14
$ws1.Range("F9").Value = 252

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 83
$ws2.Range("F3").Value = 34

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 156
$ws4.Range("F3").Value = 83
$ws4.Range("F4").Value = 469
$ws4.Range("F5").Value = 11
$ws4.Range("F7").Value = 13
$ws4.Range("F9").Value = 14
$ws4.Range("F10").Value = 252
$ws4.Range("F11").Value = 34
